# The paragraph containing the id tag was split across three runs:
#   <id>  /  p045r_1  /  </id>
# (the middle run carries different run formatting than the other two).
# Re-"finding" the full logical text and replacing it with itself causes
# Word to coalesce the matched range back into a single run, using the
# formatting of the first run in the range - which collapses the three
# runs into one run reading "<id>p045r_1</id>".
$d = $word.ActiveDocument
$d.Content.Find.Execute("<id>p045r_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p045r_1</id>", 2)
